# Rename the "Address Book" sample domain to "Loan Book" across the
# UndoRedo sequence diagram slide (docs/diagrams/UndoRedoSequenceDiagram.pptx).
#
# Four shapes on slide 1 contain the literal substring "Address" as part of
# a class/method name in the sequence diagram. We replace just the
# "Address"/"...AddressBook" portion of the affected run(s) with
# "Loan"/"...LoanBook" so everything else in each text box - other runs,
# colors, sizes, paragraph breaks - is left completely untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 6 ("Rectangle 62"): ":Address" (own paragraph) / "BookParser" -> ":Loan" / "BookParser"
# ":" and "Address" share identical run formatting, so PowerPoint collapses
# them into a single run ":Loan" once the text actually changes. Route the
# edit through a scratch value first so the "Address"->"Loan" rewrite does
# not keep the leading ":" pinned to its own unchanged run.
$shpParser = $s.Shapes.Item(6)
$trParser = $shpParser.TextFrame.TextRange
$origParser = $trParser.Characters(1, 8)
$origParser.Text = "############"
$trParser2 = $shpParser.TextFrame.TextRange
$trParser2.Characters(1, 12).Text = ":Loan"

# Shape 19 ("TextBox 78"): "undo" + "AddressBook" + "()" -> "undo" + "LoanBook" + "()"
$shpUndo = $s.Shapes.Item(19)
$trUndo = $shpUndo.TextFrame.TextRange
$trUndo.Characters(5, 11).Text = "LoanBook"

# Shape 23 ("Rectangle 62"): ":" + "VersionedAddressBook" -> ":" + "VersionedLoanBook"
$shpVersioned = $s.Shapes.Item(23)
$trVersioned = $shpVersioned.TextFrame.TextRange
$trVersioned.Characters(2, 20).Text = "VersionedLoanBook"

# Shape 35 ("TextBox 87"): "resetData" + "(" + "ReadOnlyAddressBook" + ")" -> ... + "ReadOnlyLoanBook" + ")"
$shpReadOnly = $s.Shapes.Item(35)
$trReadOnly = $shpReadOnly.TextFrame.TextRange
$trReadOnly.Characters(11, 19).Text = "ReadOnlyLoanBook"
